# invitee remark field added and remove validation of mobile no in speaker
#
# This workbook's first sheet ("invitee_sample") lists invitee import
# columns in row 1 (A1:P1) with a sample row below (row 2). The edit adds
# a new "Remark" column immediately after the existing "Profile picture"
# column (P), i.e. a new header in Q1, and leaves the sample data row
# untouched for that new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Remark" header for the invitee import template, right
# after "Profile picture" (column P).
$ws.Range("Q1").Value = "Remark"

# Mirror the saved selection state (Excel ends up with Q1 selected after
# typing the new header).
$ws.Range("Q1").Select()
